$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 4162
$wsExpo.Range("F4").Value = 2389
$wsExpo.Range("F9").Value = 209
$wsExpo.Range("F10").Value = 121
$wsExpo.Range("F11").Value = 111
$wsExpo.Range("F13").Value = 1558
$wsExpo.Range("F14").Value = 286
$wsExpo.Range("F15").Value = 3128

# Sheet "全部类型" (All types) - sheet4
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 4162
$wsAll.Range("F4").Value = 2389
$wsAll.Range("F11").Value = 209
$wsAll.Range("F12").Value = 121
$wsAll.Range("F13").Value = 111
$wsAll.Range("F17").Value = 1558
$wsAll.Range("F18").Value = 286
$wsAll.Range("F19").Value = 3128
